# Applies the text replacements described in the commit diff using
# Word COM Find/Replace. Order is chosen (via topological sort on
# substring containment) so that no replacement accidentally matches
# text freshly introduced by an earlier replacement.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("68-14=", $true, $false, $false, $false, $false, $true, 1, $false, "17+79=", 2) | Out-Null
$d.Content.Find.Execute("86-9=", $true, $false, $false, $false, $false, $true, 1, $false, "5+5=", 2) | Out-Null
$d.Content.Find.Execute("62-48=", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=", 2) | Out-Null
$d.Content.Find.Execute("97-25=", $true, $false, $false, $false, $false, $true, 1, $false, "59-32=", 2) | Out-Null
$d.Content.Find.Execute("25+10=", $true, $false, $false, $false, $false, $true, 1, $false, "98-87=", 2) | Out-Null
$d.Content.Find.Execute("86-20=", $true, $false, $false, $false, $false, $true, 1, $false, "71+4=", 2) | Out-Null
$d.Content.Find.Execute("6+44=", $true, $false, $false, $false, $false, $true, 1, $false, "38-6=", 2) | Out-Null
$d.Content.Find.Execute("18+17=", $true, $false, $false, $false, $false, $true, 1, $false, "50+14=", 2) | Out-Null
$d.Content.Find.Execute("34-16=", $true, $false, $false, $false, $false, $true, 1, $false, "88-3=", 2) | Out-Null
$d.Content.Find.Execute("49-10=", $true, $false, $false, $false, $false, $true, 1, $false, "80-34=", 2) | Out-Null
$d.Content.Find.Execute("55-42=", $true, $false, $false, $false, $false, $true, 1, $false, "98-30=", 2) | Out-Null
$d.Content.Find.Execute("67-44=", $true, $false, $false, $false, $false, $true, 1, $false, "15+39=", 2) | Out-Null
$d.Content.Find.Execute("63-7=", $true, $false, $false, $false, $false, $true, 1, $false, "3+92=", 2) | Out-Null
$d.Content.Find.Execute("71+20=", $true, $false, $false, $false, $false, $true, 1, $false, "94-72=", 2) | Out-Null
$d.Content.Find.Execute("51-39=", $true, $false, $false, $false, $false, $true, 1, $false, "51+5=", 2) | Out-Null
$d.Content.Find.Execute("12+48=", $true, $false, $false, $false, $false, $true, 1, $false, "43+7=", 2) | Out-Null
$d.Content.Find.Execute("27-15=", $true, $false, $false, $false, $false, $true, 1, $false, "64+30=", 2) | Out-Null
$d.Content.Find.Execute("57+8=", $true, $false, $false, $false, $false, $true, 1, $false, "49-1=", 2) | Out-Null
$d.Content.Find.Execute("85-41=", $true, $false, $false, $false, $false, $true, 1, $false, "27+0=", 2) | Out-Null
$d.Content.Find.Execute("79-41=", $true, $false, $false, $false, $false, $true, 1, $false, "96-27=", 2) | Out-Null
$d.Content.Find.Execute("1+44=", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=", 2) | Out-Null
$d.Content.Find.Execute("58-13=", $true, $false, $false, $false, $false, $true, 1, $false, "52+14=", 2) | Out-Null
$d.Content.Find.Execute("15+72=", $true, $false, $false, $false, $false, $true, 1, $false, "62-16=", 2) | Out-Null
$d.Content.Find.Execute("48+30=", $true, $false, $false, $false, $false, $true, 1, $false, "92-74=", 2) | Out-Null
$d.Content.Find.Execute("41-3=", $true, $false, $false, $false, $false, $true, 1, $false, "46-18=", 2) | Out-Null
$d.Content.Find.Execute("66-50=", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=", 2) | Out-Null
$d.Content.Find.Execute("88-71=", $true, $false, $false, $false, $false, $true, 1, $false, "99-52=", 2) | Out-Null
$d.Content.Find.Execute("46+45=", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=", 2) | Out-Null
$d.Content.Find.Execute("48-44=", $true, $false, $false, $false, $false, $true, 1, $false, "68-62=", 2) | Out-Null
$d.Content.Find.Execute("12+37=", $true, $false, $false, $false, $false, $true, 1, $false, "23-1=", 2) | Out-Null
$d.Content.Find.Execute("6-6=", $true, $false, $false, $false, $false, $true, 1, $false, "91+0=", 2) | Out-Null
$d.Content.Find.Execute("21-4=", $true, $false, $false, $false, $false, $true, 1, $false, "96-6=", 2) | Out-Null
$d.Content.Find.Execute("16+70=", $true, $false, $false, $false, $false, $true, 1, $false, "97-49=", 2) | Out-Null
$d.Content.Find.Execute("66+22=", $true, $false, $false, $false, $false, $true, 1, $false, "33+62=", 2) | Out-Null
$d.Content.Find.Execute("67+26=", $true, $false, $false, $false, $false, $true, 1, $false, "12+50=", 2) | Out-Null
$d.Content.Find.Execute("22+76=", $true, $false, $false, $false, $false, $true, 1, $false, "14+1=", 2) | Out-Null
$d.Content.Find.Execute("82-55=", $true, $false, $false, $false, $false, $true, 1, $false, "57-36=", 2) | Out-Null
$d.Content.Find.Execute("24-11=", $true, $false, $false, $false, $false, $true, 1, $false, "42-15=", 2) | Out-Null
$d.Content.Find.Execute("56+4=", $true, $false, $false, $false, $false, $true, 1, $false, "55+16=", 2) | Out-Null
$d.Content.Find.Execute("79-52=", $true, $false, $false, $false, $false, $true, 1, $false, "50-16=", 2) | Out-Null
$d.Content.Find.Execute("38+26=", $true, $false, $false, $false, $false, $true, 1, $false, "17-4=", 2) | Out-Null
$d.Content.Find.Execute("90-25=", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=", 2) | Out-Null
$d.Content.Find.Execute("33-21=", $true, $false, $false, $false, $false, $true, 1, $false, "27-17=", 2) | Out-Null
$d.Content.Find.Execute("57+12=", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=", 2) | Out-Null
$d.Content.Find.Execute("64-12=", $true, $false, $false, $false, $false, $true, 1, $false, "54+7=", 2) | Out-Null
$d.Content.Find.Execute("66-32=", $true, $false, $false, $false, $false, $true, 1, $false, "60+17=", 2) | Out-Null
$d.Content.Find.Execute("97-41=", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=", 2) | Out-Null
$d.Content.Find.Execute("44-8=", $true, $false, $false, $false, $false, $true, 1, $false, "92-85=", 2) | Out-Null
$d.Content.Find.Execute("37+33=", $true, $false, $false, $false, $false, $true, 1, $false, "6+6=", 2) | Out-Null
$d.Content.Find.Execute("5+75=", $true, $false, $false, $false, $false, $true, 1, $false, "35-12=", 2) | Out-Null
$d.Content.Find.Execute("94-79=", $true, $false, $false, $false, $false, $true, 1, $false, "9+77=", 2) | Out-Null
$d.Content.Find.Execute("0+78=", $true, $false, $false, $false, $false, $true, 1, $false, "73-61=", 2) | Out-Null
$d.Content.Find.Execute("47+36=", $true, $false, $false, $false, $false, $true, 1, $false, "55+37=", 2) | Out-Null
$d.Content.Find.Execute("45+40=", $true, $false, $false, $false, $false, $true, 1, $false, "98-72=", 2) | Out-Null
$d.Content.Find.Execute("20-16=", $true, $false, $false, $false, $false, $true, 1, $false, "78+11=", 2) | Out-Null
$d.Content.Find.Execute("42-31=", $true, $false, $false, $false, $false, $true, 1, $false, "11+57=", 2) | Out-Null
$d.Content.Find.Execute("40+11=", $true, $false, $false, $false, $false, $true, 1, $false, "78-20=", 2) | Out-Null
$d.Content.Find.Execute("92+5=", $true, $false, $false, $false, $false, $true, 1, $false, "1+9=", 2) | Out-Null
$d.Content.Find.Execute("21+41=", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=", 2) | Out-Null
$d.Content.Find.Execute("94-11=", $true, $false, $false, $false, $false, $true, 1, $false, "96-22=", 2) | Out-Null
$d.Content.Find.Execute("13-2=", $true, $false, $false, $false, $false, $true, 1, $false, "93+6=", 2) | Out-Null
$d.Content.Find.Execute("0+49=", $true, $false, $false, $false, $false, $true, 1, $false, "54+29=", 2) | Out-Null
$d.Content.Find.Execute("12-12=", $true, $false, $false, $false, $false, $true, 1, $false, "64+33=", 2) | Out-Null
$d.Content.Find.Execute("86-26=", $true, $false, $false, $false, $false, $true, 1, $false, "10-4=", 2) | Out-Null
$d.Content.Find.Execute("64+29=", $true, $false, $false, $false, $false, $true, 1, $false, "11-3=", 2) | Out-Null
$d.Content.Find.Execute("1+86=", $true, $false, $false, $false, $false, $true, 1, $false, "68-15=", 2) | Out-Null
$d.Content.Find.Execute("46-26=", $true, $false, $false, $false, $false, $true, 1, $false, "12+33=", 2) | Out-Null
$d.Content.Find.Execute("22+21=", $true, $false, $false, $false, $false, $true, 1, $false, "40+43=", 2) | Out-Null
$d.Content.Find.Execute("95-20=", $true, $false, $false, $false, $false, $true, 1, $false, "18+12=", 2) | Out-Null
$d.Content.Find.Execute("69-18=", $true, $false, $false, $false, $false, $true, 1, $false, "68-23=", 2) | Out-Null
$d.Content.Find.Execute("4+45=", $true, $false, $false, $false, $false, $true, 1, $false, "64-31=", 2) | Out-Null
$d.Content.Find.Execute("98-58=", $true, $false, $false, $false, $false, $true, 1, $false, "59-38=", 2) | Out-Null
$d.Content.Find.Execute("36-10=", $true, $false, $false, $false, $false, $true, 1, $false, "76+15=", 2) | Out-Null
$d.Content.Find.Execute("88-44=", $true, $false, $false, $false, $false, $true, 1, $false, "87-87=", 2) | Out-Null
$d.Content.Find.Execute("72-48=", $true, $false, $false, $false, $false, $true, 1, $false, "44+41=", 2) | Out-Null
$d.Content.Find.Execute("62-51=", $true, $false, $false, $false, $false, $true, 1, $false, "6+34=", 2) | Out-Null
$d.Content.Find.Execute("46-37=", $true, $false, $false, $false, $false, $true, 1, $false, "23+40=", 2) | Out-Null
$d.Content.Find.Execute("61-20=", $true, $false, $false, $false, $false, $true, 1, $false, "27+17=", 2) | Out-Null
$d.Content.Find.Execute("87+0=", $true, $false, $false, $false, $false, $true, 1, $false, "42+2=", 2) | Out-Null
$d.Content.Find.Execute("44-19=", $true, $false, $false, $false, $false, $true, 1, $false, "22+57=", 2) | Out-Null
$d.Content.Find.Execute("42+8=", $true, $false, $false, $false, $false, $true, 1, $false, "31-0=", 2) | Out-Null
$d.Content.Find.Execute("38-29=", $true, $false, $false, $false, $false, $true, 1, $false, "77-60=", 2) | Out-Null
$d.Content.Find.Execute("32+52=", $true, $false, $false, $false, $false, $true, 1, $false, "35+53=", 2) | Out-Null
$d.Content.Find.Execute("54+18=", $true, $false, $false, $false, $false, $true, 1, $false, "3+30=", 2) | Out-Null
$d.Content.Find.Execute("76-52=", $true, $false, $false, $false, $false, $true, 1, $false, "61-52=", 2) | Out-Null
$d.Content.Find.Execute("41+32=", $true, $false, $false, $false, $false, $true, 1, $false, "14+1=", 2) | Out-Null
$d.Content.Find.Execute("48-38=", $true, $false, $false, $false, $false, $true, 1, $false, "58-33=", 2) | Out-Null
$d.Content.Find.Execute("35+52=", $true, $false, $false, $false, $false, $true, 1, $false, "23-20=", 2) | Out-Null
$d.Content.Find.Execute("14+69=", $true, $false, $false, $false, $false, $true, 1, $false, "73-25=", 2) | Out-Null
$d.Content.Find.Execute("4+77=", $true, $false, $false, $false, $false, $true, 1, $false, "46-6=", 2) | Out-Null
$d.Content.Find.Execute("44+3=", $true, $false, $false, $false, $false, $true, 1, $false, "99-23=", 2) | Out-Null
$d.Content.Find.Execute("1+53=", $true, $false, $false, $false, $false, $true, 1, $false, "78-7=", 2) | Out-Null
$d.Content.Find.Execute("72-3=", $true, $false, $false, $false, $false, $true, 1, $false, "65-59=", 2) | Out-Null
$d.Content.Find.Execute("82-33=", $true, $false, $false, $false, $false, $true, 1, $false, "89-81=", 2) | Out-Null
$d.Content.Find.Execute("39-14=", $true, $false, $false, $false, $false, $true, 1, $false, "91+1=", 2) | Out-Null
$d.Content.Find.Execute("49+37=", $true, $false, $false, $false, $false, $true, 1, $false, "64-41=", 2) | Out-Null
$d.Content.Find.Execute("23+51=", $true, $false, $false, $false, $false, $true, 1, $false, "92-78=", 2) | Out-Null
$d.Content.Find.Execute("9+62=", $true, $false, $false, $false, $false, $true, 1, $false, "48+27=", 2) | Out-Null
$d.Content.Find.Execute("59-54=", $true, $false, $false, $false, $false, $true, 1, $false, "38+36=", 2) | Out-Null
$d.Content.Find.Execute("47+7=", $true, $false, $false, $false, $false, $true, 1, $false, "82-39=", 2) | Out-Null
